$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "عيسى لطيف"
$ws.Range("B6").Value = "ابواب وشبابيك بلاستك"
$ws.Range("C6").Value = "شارع المعامل"
$ws.Range("D6").Value = 770000000

$ws.Columns.Item(4).AutoFit() | Out-Null

$ws.Range("D6").Select() | Out-Null
